$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header M1 from "%sys_c" to "%norm_c"
$ws.Range("M1").Value2 = "%norm_c"

# ymin/ymax (F,G) and etaMin/etaMax (Q,R) become FALSE (boolean) instead of -1000/1000
# for every data row (2 through 51)
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value2 = $false   # F - ymin
    $ws.Cells.Item($r, 7).Value2 = $false   # G - ymax
    $ws.Cells.Item($r, 17).Value2 = $false  # Q - etaMin
    $ws.Cells.Item($r, 18).Value2 = $false  # R - etaMax
}

# Reset the view so the selection / top-left cell is back at A1
$ws.Range("A1").Select()
